$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the weekly block (row 68), pushing the
# existing weekly groups down by 3 rows (68-96 -> 71-99).
$ws.Rows("68:70").Insert()

# New week's data (Fecha = 44508) for the three "Calidad" grades, mirroring
# the layout of the surrounding rows.
$newRows = @(
  @{
    Row = 68
    A = 8; B = "Terminal La Palmera de La Serena"; C = "Coquimbo"; D = 44508
    E = 4; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"
    K = "Cultivar IV Región"; L = "Especial"; M = 400; N = 1800; O = 1900; P = 1850
    Q = "$/kilo (en caja de 15 kilos)"; R = "Provincia de Limarí"; S = 1850; T = 1
  },
  @{
    Row = 69
    A = 8; B = "Terminal La Palmera de La Serena"; C = "Coquimbo"; D = 44508
    E = 4; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"
    K = "Cultivar IV Región"; L = "Primera"; M = 400; N = 1500; O = 1600; P = 1550
    Q = "$/kilo (en caja de 15 kilos)"; R = "Provincia de Limarí"; S = 1550; T = 1
  },
  @{
    Row = 70
    A = 8; B = "Terminal La Palmera de La Serena"; C = "Coquimbo"; D = 44508
    E = 4; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"
    K = "Cultivar IV Región"; L = "Segunda"; M = 360; N = 1300; O = 1400; P = 1350
    Q = "$/kilo (en caja de 15 kilos)"; R = "Provincia de Limarí"; S = 1350; T = 1
  }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($entry in $newRows) {
  $r = $entry.Row
  foreach ($col in $cols) {
    $addr = "$col$r"
    $ws.Range($addr).Value = $entry[$col]
  }
}
